$d = $word.ActiveDocument

# The title paragraph ("#[ewidencja-title]") is the very first paragraph
# in the document body. Update its paragraph mark run properties and the
# run's font size from 6pt (sz/szCs = 12 half-points) to 14pt (sz/szCs = 28
# half-points). We operate directly on the paragraph/range objects so that
# only this paragraph is affected, not the many other 6pt runs used
# throughout the table cells.

$titlePara = $d.Paragraphs(1)
$titleRange = $titlePara.Range

# Update the run(s) of text within the title paragraph.
$titleRange.Font.Size = 14

# Update the paragraph mark's stored run formatting (w:pPr/w:rPr) as well,
# by including the paragraph mark in the range we format.
$fullRange = $d.Range($titlePara.Range.Start, $titlePara.Range.End)
$fullRange.Font.Size = 14
